$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$data = @(
    @("2026-02-01", "16:04:25", "16:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "16:04:29", "16:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "16:04:39", "16:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "16:04:50", "16:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 74
$endRow = $startRow + $data.Count - 1

# Force the new cells to be treated as plain text so values like
# "2026-02-01" are not auto-converted into date serial numbers,
# matching the inline string cells used throughout this sheet.
$newRange = $ws.Range("A" + $startRow + ":F" + $endRow)
$newRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $data[$i][$col - 1]
    }
}

# Restore default styling so the new rows don't carry an explicit
# number-format style that the original rows don't have.
$newRange.Style = "Normal"
